$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row for account 004886366 / RENATO / 23300.88 (currently row 3)
$ws.Rows(3).Delete()

# Insert a new row for account 005654767 / DIEGO / 400 right after the
# 004455356 / MARCELO / 418.35 row (now row 19 after the deletion above,
# originally row 20 / row 19 for FLAVIA), pushing the remaining rows down.
$ws.Rows(19).Insert()
$ws.Cells.Item(19, 1).Value = "005654767"
$ws.Cells.Item(19, 2).Value = "DIEGO"
$ws.Cells.Item(19, 3).Value = 400
